# Daily attendance processing - 2026-01-31 22:02:46
# Swap the order of names in the "Recorded By" column (G) so that
# "System, dnasr281@gmail.com" becomes "dnasr281@gmail.com, System".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = 7
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
